$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy format (style s=1) from E1, then set its value/label
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F105: time_taken timestamps (no special style, matches E column data cells)
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:38:59.957070"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:38:59.957081"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:38:59.957085"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:38:59.957088"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:38:59.957090"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:38:59.957093"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:38:59.957096"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:38:59.957098"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:38:59.957101"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:38:59.957104"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:38:59.957106"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:38:59.957111"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:38:59.957113"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:38:59.957116"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:38:59.957118"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:38:59.957121"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:38:59.957124"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:38:59.957126"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:38:59.957129"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:38:59.957131"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:38:59.957134"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:38:59.957137"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:38:59.957139"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:38:59.957142"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:38:59.957145"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:38:59.957148"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:38:59.957151"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:38:59.957153"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:38:59.957156"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:38:59.957158"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:38:59.957161"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:38:59.957163"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:38:59.957166"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:38:59.957169"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:38:59.957172"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:38:59.957175"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:38:59.957177"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:38:59.957180"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:38:59.957182"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:38:59.957185"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:38:59.957188"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:38:59.957191"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:38:59.957193"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:38:59.957196"
$ws.Cells.Item(46, 6).Value = "2021-10-05 13:38:59.957199"
$ws.Cells.Item(47, 6).Value = "2021-10-05 13:38:59.957201"
$ws.Cells.Item(48, 6).Value = "2021-10-05 13:38:59.957204"
$ws.Cells.Item(49, 6).Value = "2021-10-05 13:38:59.957206"
$ws.Cells.Item(50, 6).Value = "2021-10-05 13:38:59.957209"
$ws.Cells.Item(51, 6).Value = "2021-10-05 13:38:59.957211"
$ws.Cells.Item(52, 6).Value = "2021-10-05 13:38:59.957214"
$ws.Cells.Item(53, 6).Value = "2021-10-05 13:38:59.957216"
$ws.Cells.Item(54, 6).Value = "2021-10-05 13:38:59.957219"
$ws.Cells.Item(55, 6).Value = "2021-10-05 13:38:59.957222"
$ws.Cells.Item(56, 6).Value = "2021-10-05 13:38:59.957224"
$ws.Cells.Item(57, 6).Value = "2021-10-05 13:38:59.957227"
$ws.Cells.Item(58, 6).Value = "2021-10-05 13:38:59.957230"
$ws.Cells.Item(59, 6).Value = "2021-10-05 13:38:59.957232"
$ws.Cells.Item(60, 6).Value = "2021-10-05 13:38:59.957235"
$ws.Cells.Item(61, 6).Value = "2021-10-05 13:38:59.957237"
$ws.Cells.Item(62, 6).Value = "2021-10-05 13:38:59.957240"
$ws.Cells.Item(63, 6).Value = "2021-10-05 13:38:59.957243"
$ws.Cells.Item(64, 6).Value = "2021-10-05 13:38:59.957245"
$ws.Cells.Item(65, 6).Value = "2021-10-05 13:38:59.957248"
$ws.Cells.Item(66, 6).Value = "2021-10-05 13:38:59.957252"
$ws.Cells.Item(67, 6).Value = "2021-10-05 13:38:59.957254"
$ws.Cells.Item(68, 6).Value = "2021-10-05 13:38:59.957257"
$ws.Cells.Item(69, 6).Value = "2021-10-05 13:38:59.957260"
$ws.Cells.Item(70, 6).Value = "2021-10-05 13:38:59.957287"
$ws.Cells.Item(71, 6).Value = "2021-10-05 13:38:59.957292"
$ws.Cells.Item(72, 6).Value = "2021-10-05 13:38:59.957295"
$ws.Cells.Item(73, 6).Value = "2021-10-05 13:38:59.957297"
$ws.Cells.Item(74, 6).Value = "2021-10-05 13:38:59.957300"
$ws.Cells.Item(75, 6).Value = "2021-10-05 13:38:59.957303"
$ws.Cells.Item(76, 6).Value = "2021-10-05 13:38:59.957305"
$ws.Cells.Item(77, 6).Value = "2021-10-05 13:38:59.957308"
$ws.Cells.Item(78, 6).Value = "2021-10-05 13:38:59.957313"
$ws.Cells.Item(79, 6).Value = "2021-10-05 13:38:59.957316"
$ws.Cells.Item(80, 6).Value = "2021-10-05 13:38:59.957319"
$ws.Cells.Item(81, 6).Value = "2021-10-05 13:38:59.957322"
$ws.Cells.Item(82, 6).Value = "2021-10-05 13:38:59.957324"
$ws.Cells.Item(83, 6).Value = "2021-10-05 13:38:59.957327"
$ws.Cells.Item(84, 6).Value = "2021-10-05 13:38:59.957330"
$ws.Cells.Item(85, 6).Value = "2021-10-05 13:38:59.957332"
$ws.Cells.Item(86, 6).Value = "2021-10-05 13:38:59.957335"
$ws.Cells.Item(87, 6).Value = "2021-10-05 13:38:59.957338"
$ws.Cells.Item(88, 6).Value = "2021-10-05 13:38:59.957340"
$ws.Cells.Item(89, 6).Value = "2021-10-05 13:38:59.957343"
$ws.Cells.Item(90, 6).Value = "2021-10-05 13:38:59.957346"
$ws.Cells.Item(91, 6).Value = "2021-10-05 13:38:59.957349"
$ws.Cells.Item(92, 6).Value = "2021-10-05 13:38:59.957351"
$ws.Cells.Item(93, 6).Value = "2021-10-05 13:38:59.957354"
$ws.Cells.Item(94, 6).Value = "2021-10-05 13:38:59.957358"
$ws.Cells.Item(95, 6).Value = "2021-10-05 13:38:59.957361"
$ws.Cells.Item(96, 6).Value = "2021-10-05 13:38:59.957364"
$ws.Cells.Item(97, 6).Value = "2021-10-05 13:38:59.957367"
$ws.Cells.Item(98, 6).Value = "2021-10-05 13:38:59.957369"
$ws.Cells.Item(99, 6).Value = "2021-10-05 13:38:59.957372"
$ws.Cells.Item(100, 6).Value = "2021-10-05 13:38:59.957375"
$ws.Cells.Item(101, 6).Value = "2021-10-05 13:38:59.957378"
$ws.Cells.Item(102, 6).Value = "2021-10-05 13:38:59.957381"
$ws.Cells.Item(103, 6).Value = "2021-10-05 13:38:59.957383"
$ws.Cells.Item(104, 6).Value = "2021-10-05 13:38:59.957386"
$ws.Cells.Item(105, 6).Value = "2021-10-05 13:38:59.957389"

Write-Host "Done"
